$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextValue "D2" "56.844.03"
Set-TextValue "E2" "  -3.67%  "
Set-TextValue "D3" "2.535.51"
Set-TextValue "E3" "  -4.86%  "
Set-TextValue "E4" "  -0.02%  "
Set-TextValue "D5" "514.63"
Set-TextValue "E5" "  -1.67%  "
Set-TextValue "D6" "141.03"
Set-TextValue "E6" "  -2.22%  "
Set-TextValue "D7" "0.999"
Set-TextValue "E7" "  -0.03%  "
Set-TextValue "D8" "0.559"
Set-TextValue "E8" "  -1.79%  "
Set-TextValue "D9" "2.541.13"
Set-TextValue "E9" "  -4.91%  "
Set-TextValue "D10" "6.52"
Set-TextValue "E10" "  -6.40%  "
Set-TextValue "D11" "0.0991"
Set-TextValue "E11" "  -3.61%  "
Set-TextValue "D12" "0.322"
Set-TextValue "E12" "  -3.88%  "
Set-TextValue "E13" "  -0.36%  "
Set-TextValue "D14" "2.982.42"
Set-TextValue "E14" "  -4.78%  "
Set-TextValue "D15" "56.869.08"
Set-TextValue "E15" "  -3.61%  "
Set-TextValue "D16" "20.01"
Set-TextValue "E16" "  -4.83%  "
Set-TextValue "E17" "  -3.08%  "
Set-TextValue "D18" "2.520.65"
Set-TextValue "E18" "  -5.60%  "
Set-TextValue "D19" "330.27"
Set-TextValue "E19" "  -2.45%  "
Set-TextValue "D20" "4.26"
Set-TextValue "E20" "  -3.20%  "
Set-TextValue "D21" "10.06"
Set-TextValue "E21" "  -3.05%  "
Set-TextValue "D22" "6.13"
Set-TextValue "E22" "  -4.20%  "
Set-TextValue "D23" "0.998"
Set-TextValue "E23" "  -0.03%  "
Set-TextValue "D24" "65.17"
Set-TextValue "E24" "  +1.28%  "
Set-TextValue "D25" "0.167"
Set-TextValue "E25" "  +0.61%  "
Set-TextValue "D26" "0.998"
Set-TextValue "E26" "  +0.03%  "
Set-TextValue "D27" "0.399"
Set-TextValue "E27" "  -4.79%  "
Set-TextValue "D28" "2.651.52"
Set-TextValue "E28" "  -4.80%  "
Set-TextValue "D29" "6.90"
Set-TextValue "E29" "  -3.20%  "
Set-TextValue "E30" "  -7.11%  "
Set-TextValue "E31" "  +0.03%  "
Set-TextValue "D32" "6.23"
Set-TextValue "E32" "  -6.81%  "
Set-TextValue "D33" "1.55"
Set-TextValue "E33" "  -2.58%  "
Set-TextValue "D34" "18.47"
Set-TextValue "E34" "  -2.15%  "
Set-TextValue "D35" "148.14"
Set-TextValue "E35" "  -1.60%  "
Set-TextValue "D36" "3.98"
Set-TextValue "E36" "  -4.02%  "
Set-TextValue "D37" "1.13"
Set-TextValue "E37" "  -4.62%  "
Set-TextValue "D38" "0.837"
Set-TextValue "E38" "  -7.38%  "
Set-TextValue "D39" "35.64"
Set-TextValue "E39" "  -3.58%  "
Set-TextValue "D40" "0.819"
Set-TextValue "E40" "  -5.86%  "
Set-TextValue "D41" "1.42"
Set-TextValue "E41" "  -2.84%  "
Set-TextValue "D42" "0.999"
Set-TextValue "E42" "  -0.01%  "
Set-TextValue "D43" "3.46"
Set-TextValue "E43" "  -3.24%  "
Set-TextValue "D44" "0.0951"
Set-TextValue "E44" "  -2.06%  "
Set-TextValue "D45" "10.59"
Set-TextValue "E45" "  -0.62%  "
Set-TextValue "D46" "263.41"
Set-TextValue "E46" "  -4.35%  "
Set-TextValue "E47" "  -6.60%  "
Set-TextValue "D48" "18.58"
Set-TextValue "E48" "  -5.93%  "
Set-TextValue "D49" "0.0516"
Set-TextValue "E49" "  -3.10%  "
Set-TextValue "D50" "1.951.02"
Set-TextValue "E50" "  -4.83%  "
Set-TextValue "D51" "0.0220"
Set-TextValue "E51" "  -3.80%  "
